$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-36 down to 19-37.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new data record.
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(18, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(18, 4).Value = 44453
$ws.Cells.Item(18, 4).NumberFormat = $ws.Cells.Item(19, 4).NumberFormat
$ws.Cells.Item(18, 5).Value = 15
$ws.Cells.Item(18, 6).Value = 100112028
$ws.Cells.Item(18, 7).Value = "Sandia"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Tercera"
$ws.Cells.Item(18, 10).Value = 700
$ws.Cells.Item(18, 11).Value = 800
$ws.Cells.Item(18, 12).Value = 850
$ws.Cells.Item(18, 13).Value = 825
$ws.Cells.Item(18, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 825
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = "Hortaliza"
